# Auto-generated edit script: update NATMI TPM-derived metrics for the Cd34-Sell sheet
# (commit: "update scripts wuth new tpm")
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 212.9195706666667
$ws.Range("H2").Value = 638.7587120000001
$ws.Range("I2").Value = 0.7081219821749344
$ws.Range("J2").Value = 0.7081219821749345
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.926994
$ws.Range("N2").Value = 2.780982
$ws.Range("O2").Value = 0.002566096653125693
$ws.Range("P2").Value = 0.002566096653125693
$ws.Range("Q2").Value = 197.375164490576
$ws.Range("R2").Value = 1776.376480415184
$ws.Range("S2").Value = 0.001817109448463831
$ws.Range("T2").Value = 0.001817109448463831
$ws.Range("G3").Value = 212.9195706666667
$ws.Range("H3").Value = 638.7587120000001
$ws.Range("I3").Value = 0.7081219821749344
$ws.Range("J3").Value = 0.7081219821749345
$ws.Range("M3").Value = 93.12610233333334
$ws.Range("N3").Value = 279.378307
$ws.Range("O3").Value = 0.2577908589658698
$ws.Range("P3").Value = 0.2577908589658698
$ws.Range("Q3").Value = 19828.3697266734
$ws.Range("R3").Value = 178455.3275400606
$ws.Range("S3").Value = 0.1825473740374907
$ws.Range("T3").Value = 0.1825473740374907
$ws.Range("G4").Value = 212.9195706666667
$ws.Range("H4").Value = 638.7587120000001
$ws.Range("I4").Value = 0.7081219821749344
$ws.Range("J4").Value = 0.7081219821749345
$ws.Range("M4").Value = 264.9957936666667
$ws.Range("N4").Value = 794.9873809999999
$ws.Range("O4").Value = 0.733559029746061
$ws.Range("P4").Value = 0.733559029746061
$ws.Range("Q4").Value = 56422.79061597925
$ws.Range("R4").Value = 507805.1155438133
$ws.Range("S4").Value = 0.5194492741861024
$ws.Range("T4").Value = 0.5194492741861024
$ws.Range("G5").Value = 212.9195706666667
$ws.Range("H5").Value = 638.7587120000001
$ws.Range("I5").Value = 0.7081219821749344
$ws.Range("J5").Value = 0.7081219821749345
$ws.Range("M5").Value = 2.197830333333333
$ws.Range("N5").Value = 6.593490999999999
$ws.Range("O5").Value = 0.006084014634943477
$ws.Range("P5").Value = 0.006084014634943477
$ws.Range("Q5").Value = 467.9610909715102
$ws.Range("R5").Value = 4211.649818743592
$ws.Range("S5").Value = 0.004308224502877485
$ws.Range("T5").Value = 0.004308224502877486
$ws.Range("I6").Value = 0.2395302541481037
$ws.Range("J6").Value = 0.2395302541481037
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.926994
$ws.Range("N6").Value = 2.780982
$ws.Range("O6").Value = 0.002566096653125693
$ws.Range("P6").Value = 0.002566096653125693
$ws.Range("Q6").Value = 66.76437747031
$ws.Range("R6").Value = 600.87939723279
$ws.Range("S6").Value = 0.0006146577834917955
$ws.Range("T6").Value = 0.0006146577834917956
$ws.Range("I7").Value = 0.2395302541481037
$ws.Range("J7").Value = 0.2395302541481037
$ws.Range("M7").Value = 93.12610233333334
$ws.Range("N7").Value = 279.378307
$ws.Range("O7").Value = 0.2577908589658698
$ws.Range("P7").Value = 0.2577908589658698
$ws.Range("Q7").Value = 6707.169893787212
$ws.Range("R7").Value = 60364.52904408491
$ws.Range("S7").Value = 0.06174870996515275
$ws.Range("T7").Value = 0.06174870996515276
$ws.Range("I8").Value = 0.2395302541481037
$ws.Range("J8").Value = 0.2395302541481037
$ws.Range("M8").Value = 264.9957936666667
$ws.Range("N8").Value = 794.9873809999999
$ws.Range("O8").Value = 0.733559029746061
$ws.Range("P8").Value = 0.733559029746061
$ws.Range("Q8").Value = 19085.64585790816
$ws.Range("R8").Value = 171770.8127211734
$ws.Range("S8").Value = 0.1757095808277103
$ws.Range("T8").Value = 0.1757095808277104
$ws.Range("I9").Value = 0.2395302541481037
$ws.Range("J9").Value = 0.2395302541481037
$ws.Range("M9").Value = 2.197830333333333
$ws.Range("N9").Value = 6.593490999999999
$ws.Range("O9").Value = 0.006084014634943477
$ws.Range("P9").Value = 0.006084014634943477
$ws.Range("Q9").Value = 158.2931216279327
$ws.Range("R9").Value = 1424.638094651395
$ws.Range("S9").Value = 0.001457305571748793
$ws.Range("T9").Value = 0.001457305571748793
$ws.Range("G10").Value = 0.5110083333333334
$ws.Range("H10").Value = 1.533025
$ws.Range("I10").Value = 0.001699497292686207
$ws.Range("J10").Value = 0.001699497292686207
$ws.Range("K10").Value = 2
$ws.Range("L10").Value = 0.6666666666666666
$ws.Range("M10").Value = 0.926994
$ws.Range("N10").Value = 2.780982
$ws.Range("O10").Value = 0.002566096653125693
$ws.Range("P10").Value = 0.002566096653125693
$ws.Range("Q10").Value = 0.47370165895
$ws.Range("R10").Value = 4.26331493055
$ws.Range("S10").Value = 0.000004361074314758253
$ws.Range("T10").Value = 0.000004361074314758253
$ws.Range("G11").Value = 0.5110083333333334
$ws.Range("H11").Value = 1.533025
$ws.Range("I11").Value = 0.001699497292686207
$ws.Range("J11").Value = 0.001699497292686207
$ws.Range("M11").Value = 93.12610233333334
$ws.Range("N11").Value = 279.378307
$ws.Range("O11").Value = 0.2577908589658698
$ws.Range("P11").Value = 0.2577908589658698
$ws.Range("Q11").Value = 47.58821434318612
$ws.Range("R11").Value = 428.293929088675
$ws.Range("S11").Value = 0.0004381148668917475
$ws.Range("T11").Value = 0.0004381148668917476
$ws.Range("G12").Value = 0.5110083333333334
$ws.Range("H12").Value = 1.533025
$ws.Range("I12").Value = 0.001699497292686207
$ws.Range("J12").Value = 0.001699497292686207
$ws.Range("M12").Value = 264.9957936666667
$ws.Range("N12").Value = 794.9873809999999
$ws.Range("O12").Value = 0.733559029746061
$ws.Range("P12").Value = 0.733559029746061
$ws.Range("Q12").Value = 135.4150588619472
$ws.Range("R12").Value = 1218.735529757525
$ws.Range("S12").Value = 0.001246681585078952
$ws.Range("T12").Value = 0.001246681585078952
$ws.Range("G13").Value = 0.5110083333333334
$ws.Range("H13").Value = 1.533025
$ws.Range("I13").Value = 0.001699497292686207
$ws.Range("J13").Value = 0.001699497292686207
$ws.Range("M13").Value = 2.197830333333333
$ws.Range("N13").Value = 6.593490999999999
$ws.Range("O13").Value = 0.006084014634943477
$ws.Range("P13").Value = 0.006084014634943477
$ws.Range("Q13").Value = 1.123109615586111
$ws.Range("R13").Value = 10.107986540275
$ws.Range("S13").Value = 0.0000103397664007497
$ws.Range("T13").Value = 0.0000103397664007497
$ws.Range("G14").Value = 12.96440966666667
$ws.Range("H14").Value = 38.893229
$ws.Range("I14").Value = 0.04311667284572963
$ws.Range("J14").Value = 0.04311667284572964
$ws.Range("K14").Value = 2
$ws.Range("L14").Value = 0.6666666666666666
$ws.Range("M14").Value = 0.926994
$ws.Range("N14").Value = 2.780982
$ws.Range("O14").Value = 0.002566096653125693
$ws.Range("P14").Value = 0.002566096653125693
$ws.Range("Q14").Value = 12.017929974542
$ws.Range("R14").Value = 108.161369770878
$ws.Range("S14").Value = 0.0001106415498833423
$ws.Range("T14").Value = 0.0001106415498833423
$ws.Range("G15").Value = 12.96440966666667
$ws.Range("H15").Value = 38.893229
$ws.Range("I15").Value = 0.04311667284572963
$ws.Range("J15").Value = 0.04311667284572964
$ws.Range("M15").Value = 93.12610233333334
$ws.Range("N15").Value = 279.378307
$ws.Range("O15").Value = 0.2577908589658698
$ws.Range("P15").Value = 0.2577908589658698
$ws.Range("Q15").Value = 1207.324941309256
$ws.Range("R15").Value = 10865.9244717833
$ws.Range("S15").Value = 0.01111508412865104
$ws.Range("T15").Value = 0.01111508412865104
$ws.Range("G16").Value = 12.96440966666667
$ws.Range("H16").Value = 38.893229
$ws.Range("I16").Value = 0.04311667284572963
$ws.Range("J16").Value = 0.04311667284572964
$ws.Range("M16").Value = 264.9957936666667
$ws.Range("N16").Value = 794.9873809999999
$ws.Range("O16").Value = 0.733559029746061
$ws.Range("P16").Value = 0.733559029746061
$ws.Range("Q16").Value = 3435.514029038139
$ws.Range("R16").Value = 30919.62626134324
$ws.Range("S16").Value = 0.03162862469859177
$ws.Range("T16").Value = 0.03162862469859177
$ws.Range("G17").Value = 12.96440966666667
$ws.Range("H17").Value = 38.893229
$ws.Range("I17").Value = 0.04311667284572963
$ws.Range("J17").Value = 0.04311667284572964
$ws.Range("M17").Value = 2.197830333333333
$ws.Range("N17").Value = 6.593490999999999
$ws.Range("O17").Value = 0.006084014634943477
$ws.Range("P17").Value = 0.006084014634943477
$ws.Range("Q17").Value = 28.49357281915989
$ws.Range("R17").Value = 256.4421553724389
$ws.Range("S17").Value = 0.0002623224686034891
$ws.Range("T17").Value = 0.0002623224686034891
$ws.Range("G18").Value = 1.873532666666667
$ws.Range("H18").Value = 5.620598
$ws.Range("I18").Value = 0.00623094279889598
$ws.Range("J18").Value = 0.006230942798895981
$ws.Range("K18").Value = 2
$ws.Range("L18").Value = 0.6666666666666666
$ws.Range("M18").Value = 0.926994
$ws.Range("N18").Value = 2.780982
$ws.Range("O18").Value = 0.002566096653125693
$ws.Range("P18").Value = 0.002566096653125693
$ws.Range("Q18").Value = 1.736753540804
$ws.Range("R18").Value = 15.630781867236
$ws.Range("S18").Value = 0.00001598920146206461
$ws.Range("T18").Value = 0.00001598920146206462
$ws.Range("G19").Value = 1.873532666666667
$ws.Range("H19").Value = 5.620598
$ws.Range("I19").Value = 0.00623094279889598
$ws.Range("J19").Value = 0.006230942798895981
$ws.Range("M19").Value = 93.12610233333334
$ws.Range("N19").Value = 279.378307
$ws.Range("O19").Value = 0.2577908589658698
$ws.Range("P19").Value = 0.2577908589658698
$ws.Range("Q19").Value = 174.4747948408429
$ws.Range("R19").Value = 1570.273153567586
$ws.Range("S19").Value = 0.001606280096294596
$ws.Range("T19").Value = 0.001606280096294596
$ws.Range("G20").Value = 1.873532666666667
$ws.Range("H20").Value = 5.620598
$ws.Range("I20").Value = 0.00623094279889598
$ws.Range("J20").Value = 0.006230942798895981
$ws.Range("M20").Value = 264.9957936666667
$ws.Range("N20").Value = 794.9873809999999
$ws.Range("O20").Value = 0.733559029746061
$ws.Range("P20").Value = 0.733559029746061
$ws.Range("Q20").Value = 496.4782759637598
$ws.Range("R20").Value = 4468.304483673837
$ws.Range("S20").Value = 0.004570764353961341
$ws.Range("T20").Value = 0.004570764353961342
$ws.Range("G21").Value = 1.873532666666667
$ws.Range("H21").Value = 5.620598
$ws.Range("I21").Value = 0.00623094279889598
$ws.Range("J21").Value = 0.006230942798895981
$ws.Range("M21").Value = 2.197830333333333
$ws.Range("N21").Value = 6.593490999999999
$ws.Range("O21").Value = 0.006084014634943477
$ws.Range("P21").Value = 0.006084014634943477
$ws.Range("Q21").Value = 4.117706925290888
$ws.Range("R21").Value = 37.059362327618
$ws.Range("S21").Value = 0.00003790914717797881
$ws.Range("T21").Value = 0.00003790914717797881
$ws.Range("G22").Value = 0.3910823333333333
$ws.Range("H22").Value = 1.173247
$ws.Range("I22").Value = 0.001300650739650178
$ws.Range("J22").Value = 0.001300650739650178
$ws.Range("K22").Value = 2
$ws.Range("L22").Value = 0.6666666666666666
$ws.Range("M22").Value = 0.926994
$ws.Range("N22").Value = 2.780982
$ws.Range("O22").Value = 0.002566096653125693
$ws.Range("P22").Value = 0.002566096653125693
$ws.Range("Q22").Value = 0.362530976506
$ws.Range("R22").Value = 3.262778788554
$ws.Range("S22").Value = 0.000003337595509901779
$ws.Range("T22").Value = 0.00000333759550990178
$ws.Range("G23").Value = 0.3910823333333333
$ws.Range("H23").Value = 1.173247
$ws.Range("I23").Value = 0.001300650739650178
$ws.Range("J23").Value = 0.001300650739650178
$ws.Range("M23").Value = 93.12610233333334
$ws.Range("N23").Value = 279.378307
$ws.Range("O23").Value = 0.2577908589658698
$ws.Range("P23").Value = 0.2577908589658698
$ws.Range("Q23").Value = 36.41997339475878
$ws.Range("R23").Value = 327.779760552829
$ws.Range("S23").Value = 0.0003352958713890133
$ws.Range("T23").Value = 0.0003352958713890133
$ws.Range("G24").Value = 0.3910823333333333
$ws.Range("H24").Value = 1.173247
$ws.Range("I24").Value = 0.001300650739650178
$ws.Range("J24").Value = 0.001300650739650178
$ws.Range("M24").Value = 264.9957936666667
$ws.Range("N24").Value = 794.9873809999999
$ws.Range("O24").Value = 0.733559029746061
$ws.Range("P24").Value = 0.733559029746061
$ws.Range("Q24").Value = 103.6351733106785
$ws.Range("R24").Value = 932.7165597961068
$ws.Range("S24").Value = 0.0009541040946162813
$ws.Range("T24").Value = 0.0009541040946162814
$ws.Range("G25").Value = 0.3910823333333333
$ws.Range("H25").Value = 1.173247
$ws.Range("I25").Value = 0.001300650739650178
$ws.Range("J25").Value = 0.001300650739650178
$ws.Range("M25").Value = 2.197830333333333
$ws.Range("N25").Value = 6.593490999999999
$ws.Range("O25").Value = 0.006084014634943477
$ws.Range("P25").Value = 0.006084014634943477
$ws.Range("Q25").Value = 0.8595326150307776
$ws.Range("R25").Value = 7.735793535276999
$ws.Range("S25").Value = 0.000007913178134981741
$ws.Range("T25").Value = 0.000007913178134981743
